# Add data for 2022-11-30 (new day of carjacking data: 11-22)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet/tab name and workbook sheet name to reflect the new "through" date
$ws.Name = "Through 2022-11-22"

# Update the row label for November to reflect new "through" date
$ws.Range("A12").Value = "November (through 11-22)"

# Update November row (row 12) values for each year column (B:I)
$ws.Range("B12").Value = 22
$ws.Range("C12").Value = 53
$ws.Range("D12").Value = 88
$ws.Range("E12").Value = 43
$ws.Range("F12").Value = 36
$ws.Range("G12").Value = 152
$ws.Range("H12").Value = 151
$ws.Range("I12").Value = 83

# Update Total row (row 13) values for each year column (B:I)
$ws.Range("B13").Value = 280
$ws.Range("C13").Value = 539
$ws.Range("D13").Value = 798
$ws.Range("E13").Value = 658
$ws.Range("F13").Value = 518
$ws.Range("G13").Value = 1209
$ws.Range("H13").Value = 1592
$ws.Range("I13").Value = 1480
